$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-15 down to 5-16.
$ws.Rows("4").Insert()

# Fill in the new module row: "app.fancy-slider.blur" is incomplete.
$ws.Range("A4").Value = "app.fancy-slider.blur"
$ws.Range("B4").Value = "Nu-I complet!"
$ws.Range("D4").Value = "Trebuie integrat!"

# Update the selection to match the author's last-saved cursor position.
$ws.Range("E4").Select()
